$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = 3
    "F4"  = 0
    "F5"  = -3
    "F15" = -6
    "F21" = 3
    "F25" = -1
    "F33" = 8
    "F34" = -1
    "F36" = 0
    "F39" = 0
    "F41" = -3
    "F46" = -3
    "F49" = -2
    "F53" = -6
    "F56" = -1
    "F57" = -4
    "F59" = 0
    "F60" = -3
    "F61" = -1
    "F62" = -3
    "F64" = -2
    "F66" = 0
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
